$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.976.18"
$ws.Range("E2").Value = "  -0.84%  "

# Row 3
$ws.Range("D3").Value = "3.780.59"
$ws.Range("E3").Value = "  -1.67%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.47"
$ws.Range("E5").Value = "  -0.79%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.69"
$ws.Range("E6").Value = "  +0.78%  "

# Row 7
$ws.Range("D7").Value = "3.777.99"
$ws.Range("E7").Value = "  -1.79%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  +0.05%  "

# Row 10
$ws.Range("E10").Value = "  -1.21%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.55"
$ws.Range("E11").Value = "  +1.43%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").Value = "  -0.72%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000278"
$ws.Range("E13").Value = "  +0.52%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.39"
$ws.Range("E14").Value = "  -0.97%  "

# Row 15
$ws.Range("D15").Value = "4.409.93"
$ws.Range("E15").Value = "  -1.62%  "

# Row 16
$ws.Range("D16").Value = "3.787.26"
$ws.Range("E16").Value = "  -1.43%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.13"
$ws.Range("E17").Value = "  +4.48%  "

# Row 18
$ws.Range("D18").Value = "67.917.08"
$ws.Range("E18").Value = "  -0.77%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.21"
$ws.Range("E19").Value = "  -1.49%  "

# Row 20
$ws.Range("E20").Value = "  +0.86%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.61"
$ws.Range("E21").Value = "  -3.28%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.04"
$ws.Range("E22").Value = "  -0.64%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.721"
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000150"
$ws.Range("E24").Value = "  -7.62%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.85"
$ws.Range("E25").Value = "  +0.76%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("E26").Value = "  +0.47%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.16"
$ws.Range("E27").Value = "  +1.19%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.55"
$ws.Range("E28").Value = "  +2.03%  "

# Row 29
$ws.Range("E29").Value = "  +0.10%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.93"
$ws.Range("E30").Value = "  -0.66%  "

# Row 31
$ws.Range("D31").Value = "3.925.61"
$ws.Range("E31").Value = "  -1.59%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.61"
$ws.Range("E32").Value = "  -0.58%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.55"
$ws.Range("E33").Value = "  -1.97%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.25"
$ws.Range("E34").Value = "  -1.67%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.21"
$ws.Range("E35").Value = "  -0.26%  "

# Row 36
$ws.Range("D36").Value = "3.738.39"
$ws.Range("E36").Value = "  -1.79%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.105"
$ws.Range("E37").Value = "  +0.87%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.75"
$ws.Range("E38").Value = "  -0.51%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.140"
$ws.Range("E39").Value = "  -0.03%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.01"
$ws.Range("E40").Value = "  -1.46%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.85"
$ws.Range("E41").Value = "  -0.77%  "

# Row 42
$ws.Range("E42").Value = "  +0.09%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.314"
$ws.Range("E43").Value = "  -0.20%  "

# Row 44
$ws.Range("E44").Value = "  +0.01%  "

# Row 45
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.96"
$ws.Range("E45").Value = "  -1.52%  "

# Row 46
$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.67"
$ws.Range("E46").Value = "  +0.22%  "

# Row 47
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "403.92"
$ws.Range("E47").Value = "  -3.62%  "

# Row 48
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.88"
$ws.Range("E48").Value = "  -2.22%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000281"
$ws.Range("E49").Value = "  -4.62%  "

# Row 50
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.16"
$ws.Range("E50").Value = "  +7.65%  "

# Row 51
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "140.45"
$ws.Range("E51").Value = "  -1.22%  "
